# Refresh the cryptos price list (prices in column D, 1h volume % in column E).
# Rows 41/42 and 44/45 swap coin identities (Coin/Link/Price/Volume) per the
# updated ranking; row index in column A is unaffected.
#
# Note: several Price values (column D) are plain decimal numbers (e.g.
# "598.43", "1.00"). Assigning those bare strings via .Value lets Excel's
# COM layer auto-coerce them into numeric cells, which would not match the
# source data (stored as literal text). Prefixing with a leading apostrophe
# forces Excel to keep them as text, exactly like the original cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.441.49'
$ws.Range('E2').Value = '  -2.07%  '
$ws.Range('D3').Value = '2.641.87'
$ws.Range('E3').Value = '  -3.30%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''598.43'
$ws.Range('E5').Value = '  -0.75%  '
$ws.Range('D6').Value = '''167.12'
$ws.Range('E6').Value = '  -1.24%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -0.49%  '
$ws.Range('D9').Value = '2.641.27'
$ws.Range('E9').Value = '  -3.25%  '
$ws.Range('D10').Value = '''0.145'
$ws.Range('E10').Value = '  -0.15%  '
$ws.Range('E11').Value = '  +1.37%  '
$ws.Range('D12').Value = '''0.365'
$ws.Range('E12').Value = '  -1.17%  '
$ws.Range('D13').Value = '''5.23'
$ws.Range('E13').Value = '  -2.07%  '
$ws.Range('D14').Value = '''28.03'
$ws.Range('D15').Value = '3.124.27'
$ws.Range('E15').Value = '  -3.15%  '
$ws.Range('E16').Value = '  -3.06%  '
$ws.Range('D17').Value = '67.383.22'
$ws.Range('E17').Value = '  -1.94%  '
$ws.Range('D18').Value = '2.641.58'
$ws.Range('E18').Value = '  -2.60%  '
$ws.Range('D19').Value = '''11.92'
$ws.Range('E19').Value = '  +0.12%  '
$ws.Range('D20').Value = '''7.85'
$ws.Range('E20').Value = '  +2.45%  '
$ws.Range('D21').Value = '''363.67'
$ws.Range('E21').Value = '  -2.64%  '
$ws.Range('D22').Value = '''4.41'
$ws.Range('E22').Value = '  -3.03%  '
$ws.Range('D23').Value = '''4.79'
$ws.Range('E23').Value = '  -3.39%  '
$ws.Range('D24').Value = '''10.93'
$ws.Range('E24').Value = '  +8.80%  '
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('D27').Value = '''70.93'
$ws.Range('E27').Value = '  -3.61%  '
$ws.Range('D28').Value = '2.779.08'
$ws.Range('E28').Value = '  -3.26%  '
$ws.Range('D29').Value = '''0.0000103'
$ws.Range('E29').Value = '  -2.90%  '
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('D31').Value = '''554.39'
$ws.Range('E31').Value = '  -5.81%  '
$ws.Range('D32').Value = '''8.07'
$ws.Range('E32').Value = '  -2.81%  '
$ws.Range('D33').Value = '''1.39'
$ws.Range('E33').Value = '  -3.83%  '
$ws.Range('D34').Value = '''1.93'
$ws.Range('E34').Value = '  -1.65%  '
$ws.Range('D35').Value = '''0.133'
$ws.Range('E35').Value = '  +0.18%  '
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('E37').Value = '  -5.15%  '
$ws.Range('D38').Value = '''157.73'
$ws.Range('E38').Value = '  -2.78%  '
$ws.Range('D39').Value = '''19.42'
$ws.Range('E39').Value = '  -2.85%  '
$ws.Range('D40').Value = '''0.373'
$ws.Range('E40').Value = '  -2.41%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = '''1.83'
$ws.Range('E41').Value = '  -5.13%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').Value = '''5.27'
$ws.Range('E42').Value = '  -4.17%  '
$ws.Range('D43').Value = '''17.93'
$ws.Range('E43').Value = '  -0.39%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').Value = '''2.52'
$ws.Range('E44').Value = '  -4.97%  '
$ws.Range('B45').Value = 'USDe'
$ws.Range('C45').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D45').Value = '''1.00'
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('D46').Value = '''40.13'
$ws.Range('E46').Value = '  -2.28%  '
$ws.Range('D47').Value = '0.0₆0301'
$ws.Range('E47').Value = '  -3.38%  '
$ws.Range('D48').Value = '''0.596'
$ws.Range('E48').Value = '  -1.60%  '
$ws.Range('D49').Value = '''153.98'
$ws.Range('E49').Value = '  -1.46%  '
$ws.Range('D50').Value = '''3.89'
$ws.Range('E50').Value = '  -2.19%  '
$ws.Range('E51').Value = '  -3.75%  '
